# update matrix mult seq with one add one mult
#
# Restructure the "MatmultSeq" 1024-row (row 10) to use the same
# D:J layout (two addends, XOR+IV sum, NOR, total, weight, weighted total)
# as the other rows, instead of the old L:Q layout. Also split the old
# single "3x3" row (row 23) into three distinct rows for 2x2 / 3x3 / 5x5
# with their own multiplier (I column) and matching data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: move data from L:Q into D:J, clear L:Q ---------------------
$ws.Range("D10").Value = 2047
$ws.Range("E10").Value = 2050
$ws.Range("F10").Formula = "=SUM(D10:E10)"
$ws.Range("G10").Value = 3070
$ws.Range("H10").Formula = "=SUM(D10:F10)"
$ws.Range("I10").Value = 1
$ws.Range("J10").Formula = "=(SUM(D10:E10)+5*G10)*I10"

$ws.Range("L10").Value = $null
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = $null
$ws.Range("O10").Value = $null
$ws.Range("P10").Value = $null
$ws.Range("Q10").Value = $null

# --- Row 23: becomes the "2x2" case (was mislabeled "3x3") --------------
$ws.Range("B23").Value = "2x2"
$ws.Range("E23").Value = 997
$ws.Range("G23").Value = 1956
$ws.Range("I23").Formula = "=2*2*2"

# --- Row 24: new "3x3" row, same shape as the old row 23 ----------------
$ws.Range("B24").Value = "3x3"
$ws.Range("D24").Value = 1026
$ws.Range("E24").Value = 997
$ws.Range("F24").Formula = "=SUM(D24:E24)"
$ws.Range("G24").Value = 1956
$ws.Range("H24").Formula = "=SUM(D24:F24)"
$ws.Range("I24").Formula = "=3*3*3"
$ws.Range("J24").Formula = "=(SUM(D24:E24)+5*G24)*I24"

# Row 24 is brand new, so copy the number formats / bold weight that the
# D:J block already uses elsewhere (e.g. row 23) instead of leaving it
# defaulted to "General".
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("H24").NumberFormat = "#,##0"
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("J24").Font.Bold = $true

# --- Row 25: new "5x5" row ------------------------------------------------
$ws.Range("B25").Value = "5x5"
$ws.Range("D25").Value = 1026
$ws.Range("E25").Value = 997
$ws.Range("F25").Formula = "=SUM(D25:E25)"
$ws.Range("G25").Value = 1956
$ws.Range("H25").Formula = "=SUM(D25:F25)"
$ws.Range("I25").Formula = "=5*5*5"
$ws.Range("J25").Formula = "=(SUM(D25:E25)+5*G25)*I25"

$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("H25").NumberFormat = "#,##0"
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("J25").Font.Bold = $true

# --- View bookkeeping: scroll/selection like the authored commit --------
[void]$ws.Range("G1").Select()
$excel.ActiveWindow.ScrollColumn = 7
[void]$ws.Range("M12").Select()
